# -----------------------------------------------------------------------------
# Workbook / sheet-level changes
# -----------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook

# Rename the original (only) sheet from "Sheet" to "Results"
$results = $wb.Worksheets.Item(1)
$results.Name = "Results"

# Add a new worksheet "h2h" right after "Results"
$h2h = $wb.Worksheets.Add($null, $results)
$h2h.Name = "h2h"

# -----------------------------------------------------------------------------
# "Results" sheet --------------------------------------------------------------
# -----------------------------------------------------------------------------

# Header row translations (Romanian -> English)
$results.Range("A1").Value = "Date"
$results.Range("B1").Value = "Opponent"
$results.Range("C1").Value = "Result"
$results.Range("M1").Value = "Type"
$results.Range("O1").Value = "City"
$results.Range("P1").Value = "Venue"
$results.Range("Q1").Value = "Surface"
$results.Range("R1").Value = "Rating"
$results.Range("S1").Value = "Observations"

# Row 2 - update the existing match entry
$results.Range("A2").Value = ""
$results.Range("B2").Value = "cris"
$results.Range("C2").Value = "W"
$results.Range("D2").Value = "6-4"
$results.Range("E2").Value = "6-4"

# Row 3 - new match entry
$results.Range("B3").Value = "test"
$results.Range("C3").Value = "L"
$results.Range("D3").Value = "5-7"
$results.Range("E3").Value = "6-4"
$results.Range("F3").Value = "4-6"

# Row 4 - new match entry
$results.Range("A4").Value = "sdfd"
$results.Range("B4").Value = "cris"
$results.Range("C4").Value = "L"
$results.Range("D4").Value = "4-6"
$results.Range("E4").Value = "4-6"
$results.Range("F4").Value = "6-4"
$results.Range("G4").Value = "6-4"
$results.Range("H4").Value = "6-4"

# Row 5 - new match entry
$results.Range("B5").Value = "test"
$results.Range("C5").Value = "NA"

# -----------------------------------------------------------------------------
# "h2h" sheet --------------------------------------------------------------------
# -----------------------------------------------------------------------------

$h2h.Range("A1").Value = "Won"
$h2h.Range("C1").Value = "Lost"
$h2h.Range("E1").Value = "Opponent"

# "1" / "0" are plain text values in the source data (not numbers) - prefix
# with a quote so Excel stores them as text, same as a user typing '1.
$h2h.Range("A2").Value = "'1"
$h2h.Range("C2").Value = "'1"
$h2h.Range("E2").Value = "cris"

$h2h.Range("A3").Value = "'0"
$h2h.Range("C3").Value = "'1"
$h2h.Range("E3").Value = "test"

# Row 4 is blank in the source data, but still present - touch it so the
# sheet's used range extends down to row 4 (A4:E4).
$h2h.Range("A4:E4").Font.Bold = $false

# Activate the "h2h" sheet (second tab, i.e. activeTab = 1)
$h2h.Activate()
